$wb = $excel.ActiveWorkbook

# Add the new worksheet after the existing LoginPage sheet
$loginSheet = $wb.Worksheets.Item("LoginPage")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$newSheet.Name = "AdminUsersPage"

# Populate header row with test data (order matches shared-string insertion order)
$newSheet.Range("C1").Value = "Staff"
$newSheet.Range("A1").Value = "Simi"
$newSheet.Range("B1").Value = "Simi123"

# Auto-fit columns to mirror the committed widths
$newSheet.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$newSheet.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$newSheet.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Make the new sheet the active tab, matching the diff's activeTab/tabSelected change
$newSheet.Activate()
